$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 3: Version value 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Row 8: Date value -> new date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Row 9: Publisher value (was empty) -> Alvearie Team
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail"
# Row 11 was a duplicate "Contact" / "No display for ContactDetail"
# New layout: row 10 becomes "Jurisdiction" / "United States of America"
# and the old row 11 (duplicate) is deleted, shifting everything up.
$ws.Range("A11").EntireRow.Delete()

$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
